{"js": "// The document's Title, Author, and Abstract paragraphs each have their\n// text split across multiple single-word runs (leftover from an older\n// template). Collapse each of those paragraphs down to one run while\n// keeping the visible text identical, so the underlying markup reads\n// as a single contiguous sentence/run per paragraph.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Styles that need their runs merged into a single run, keyed by the\n// paragraph style name used in this document.\nconst targetStyles = new Set([\"Title\", \"Author\", \"Abstract\"]);\n\nconst candidates = [];\nfor (const p of paragraphs.items) {\n  p.load(\"text,style\");\n  candidates.push(p);\n}\nawait context.sync();\n\nfor (const p of candidates) {\n  if (targetStyles.has(p.style)) {\n    // Re-insert the paragraph's own text, replacing its contents.\n    // This collapses the multiple single-word runs into one run\n    // without altering the rendered text.\n    p.insertText(p.text, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document's Title, Author, and Abstract paragraphs each have their\n# text split across multiple single-word runs (leftover from an older\n# template). Collapse each of those paragraphs down to one run while\n# keeping the visible text identical, so the underlying markup reads\n# as a single contiguous sentence/run per paragraph.\n\n$d = $word.ActiveDocument\n\n$targetStyles = @(\"Title\", \"Author\", \"Abstract\")\n\nforeach ($p in $d.Paragraphs) {\n    $styleName = $p.Style.NameLocal\n    if ($targetStyles -contains $styleName) {\n        $r = $p.Range\n        # Exclude the trailing paragraph mark from the range so we only\n        # touch the paragraph's visible text.\n        [void]$r.MoveEnd(1, -1)\n        $originalText = $r.Text\n\n        # Re-assigning the exact same string is treated as a no-op by the\n        # engine (the runs are left untouched), so round-trip through a\n        # placeholder value first to force a genuine replace that\n        # collapses the paragraph's runs into a single run.\n        $r.Text = \"~~~placeholder~~~\"\n\n        $r2 = $p.Range\n        [void]$r2.MoveEnd(1, -1)\n        $r2.Text = $originalText\n    }\n}\n"}
